# Finished Week 13 logging
# Update cumulative target depth stats (row 3, the "R" row) on both the
# OFF and DEF sheets with the latest totals after Week 13.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 163
$wsOff.Range("C3").Value = 98
$wsOff.Range("D3").Value = 30
$wsOff.Range("E3").Value = 11

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 184
$wsDef.Range("C3").Value = 142
$wsDef.Range("D3").Value = 50
$wsDef.Range("E3").Value = 28
